# Update TPM-derived NATMI ligand-receptor metrics (Pdgfa-Pdgfra) to new values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2"  = "1.425879333333333";   "H2"  = "4.277638";
    "I2"  = "0.07659591414135564"; "J2"  = "0.07659591414135564";
    "K2"  = "3";                   "L2"  = "1";
    "M2"  = "0.2662156666666667";  "N2"  = "0.7986470000000001";
    "O2"  = "0.0009813702709097034"; "P2" = "0.0009813702709097034";
    "Q2"  = "0.3795914173095556";  "R2"  = "3.416322755786";
    "S2"  = "7.516895301147857E-05"; "T2" = "7.516895301147857E-05";

    "G3"  = "1.425879333333333";   "H3"  = "4.277638";
    "I3"  = "0.07659591414135564"; "J3"  = "0.07659591414135564";
    "O3"  = "0.998256289001958";   "P3"  = "0.998256289001958";
    "Q3"  = "386.1228843107029";   "R3"  = "3475.105958796326";
    "S3"  = "0.07646235300346227"; "T3"  = "0.07646235300346227";

    "G4"  = "1.425879333333333";   "H4"  = "4.277638";
    "I4"  = "0.07659591414135564"; "J4"  = "0.07659591414135564";
    "M4"  = "0.2067996666666667";  "N4"  = "0.620399";
    "O4"  = "0.000762340727132399"; "P4" = "0.0007623407271323989";
    "Q4"  = "0.2948713708402222";  "R4"  = "2.653842337562";
    "S4"  = "5.839218488189186E-05"; "T4" = "5.839218488189185E-05";

    "I5"  = "0.03783651626913671"; "J5"  = "0.03783651626913671";
    "K5"  = "3";                   "L5"  = "1";
    "M5"  = "0.2662156666666667";  "N5"  = "0.7986470000000001";
    "O5"  = "0.0009813702709097034"; "P5" = "0.0009813702709097034";
    "Q5"  = "0.1875089160781112";  "R5"  = "1.687580244703";
    "S5"  = "3.71316322213221E-05"; "T5" = "3.713163222132209E-05";

    "I6"  = "0.03783651626913671"; "J6"  = "0.03783651626913671";
    "O6"  = "0.998256289001958";   "P6"  = "0.998256289001958";
    "S6"  = "0.03777054031959062"; "T6"  = "0.03777054031959062";

    "I7"  = "0.03783651626913671"; "J7"  = "0.03783651626913671";
    "M7"  = "0.2067996666666667";  "N7"  = "0.620399";
    "O7"  = "0.000762340727132399"; "P7" = "0.0007623407271323989";
    "Q7"  = "0.1456592762834445";
    "S7"  = "2.884431732477053E-05"; "T7" = "2.884431732477052E-05";

    "G8"  = "16.48537666666666";   "H8"  = "49.45612999999999";
    "I8"  = "0.8855675695895077";  "J8"  = "0.8855675695895077";
    "K8"  = "3";                   "L8"  = "1";
    "M8"  = "0.2662156666666667";  "N8"  = "0.7986470000000001";
    "O8"  = "0.0009813702709097034"; "P8" = "0.0009813702709097034";
    "Q8"  = "4.388665539567778";   "R8"  = "39.49798985611";
    "S8"  = "0.0008690696856769027"; "T8" = "0.0008690696856769027";

    "G9"  = "16.48537666666666";   "H9"  = "49.45612999999999";
    "I9"  = "0.8855675695895077";  "J9"  = "0.8855675695895077";
    "O9"  = "0.998256289001958";   "P9"  = "0.998256289001958";
    "Q9"  = "4464.179428564334";   "R9"  = "40177.61485707901";
    "S9"  = "0.8840233956789051";  "T9"  = "0.8840233956789051";

    "G10" = "16.48537666666666";   "H10" = "49.45612999999999";
    "I10" = "0.8855675695895077";  "J10" = "0.8855675695895077";
    "M10" = "0.2067996666666667";  "N10" = "0.620399";
    "O10" = "0.000762340727132399"; "P10" = "0.0007623407271323989";
    "Q10" = "3.409170399541111";
    "S10" = "0.0006751042249257366"; "T10" = "0.0006751042249257365";
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = [double]$updates[$addr]
}

Write-Output "Applied $($updates.Count) cell updates"
